$d = $word.ActiveDocument

# --- Edit 1 ---------------------------------------------------------------
# Row "On voit tous les trajets qui nous impliquent" / "GET URLBackend/rides?mine=true"
# Response cell currently reads "Les trajets qui nous impliquent".
# Append clarifying text " (conducteur ou inscrit)".
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("Les trajets qui nous impliquent", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $rng1.Collapse(0)
    $rng1.InsertAfter(" (conducteur ou inscrit)")
}

# --- Edit 2 ---------------------------------------------------------------
# Row "On voit tous les trajets" / "GET URLBackend/rides"
# Response cell currently reads "Tous les trajets (combien ?)".
# Append clarifying text " sauf ceux qu'on conduit".
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("(combien", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $rng2.Collapse(0)
    $rng2.MoveEnd(1, 3) | Out-Null
    $rng2.Collapse(0)
    $rng2.InsertAfter(" sauf ceux qu’on conduit")
}
